# Update TPM-derived NATMI metrics (ligand/receptor expression + specificity
# scores) for the Rtn4-Rtn4r sheet, rows 2-6, per the "update scripts wuth new tpm"
# commit. Only the affected numeric cells are rewritten; everything else is
# left untouched.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value2 = 48.46865866666667
$ws.Range("H2").Value2 = 145.405976
$ws.Range("I2").Value2 = 0.1554430998624896
$ws.Range("J2").Value2 = 0.1554430998624896
$ws.Range("M2").Value2 = 0.073876
$ws.Range("Q2").Value2 = 3.580670627658667
$ws.Range("R2").Value2 = 32.226035648928
$ws.Range("S2").Value2 = 0.1554430998624896
$ws.Range("T2").Value2 = 0.1554430998624896
$ws.Range("I3").Value2 = 0.20693808715897
$ws.Range("J3").Value2 = 0.20693808715897
$ws.Range("M3").Value2 = 0.073876
$ws.Range("Q3").Value2 = 4.766870521042666
$ws.Range("S3").Value2 = 0.20693808715897
$ws.Range("T3").Value2 = 0.20693808715897
$ws.Range("G4").Value2 = 75.47903666666667
$ws.Range("H4").Value2 = 226.43711
$ws.Range("I4").Value2 = 0.2420676733554854
$ws.Range("J4").Value2 = 0.2420676733554854
$ws.Range("M4").Value2 = 0.073876
$ws.Range("Q4").Value2 = 5.576089312786667
$ws.Range("R4").Value2 = 50.18480381508
$ws.Range("S4").Value2 = 0.2420676733554854
$ws.Range("T4").Value2 = 0.2420676733554854
$ws.Range("G5").Value2 = 51.18999233333333
$ws.Range("H5").Value2 = 153.569977
$ws.Range("I5").Value2 = 0.164170647777855
$ws.Range("J5").Value2 = 0.164170647777855
$ws.Range("M5").Value2 = 0.073876
$ws.Range("Q5").Value2 = 3.781711873617333
$ws.Range("R5").Value2 = 34.03540686255599
$ws.Range("S5").Value2 = 0.164170647777855
$ws.Range("T5").Value2 = 0.164170647777855
$ws.Range("G6").Value2 = 72.14667033333333
$ws.Range("H6").Value2 = 216.440011
$ws.Range("I6").Value2 = 0.2313804918452
$ws.Range("J6").Value2 = 0.2313804918452
$ws.Range("M6").Value2 = 0.073876
$ws.Range("Q6").Value2 = 5.329907417545333
$ws.Range("R6").Value2 = 47.969166757908
$ws.Range("S6").Value2 = 0.2313804918452
$ws.Range("T6").Value2 = 0.2313804918452
